$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------- Row 8 ----------
$ws.Range("A8").Value = "R. mtr."
$ws.Range("C8").Value = 27
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "17"
$ws.Range("E8").Value = "25 mm"
$ws.Range("F8").Value = 56
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "1512.00"

# ---------- Row 9 ----------
$ws.Range("A9").Value = "Set"
$ws.Range("C9").Value = 48
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "13.0"
$e9 = @'
Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   
'@
$ws.Range("E9").Value = $e9
$ws.Range("F9").Value = 5733
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "275184.00"

# ---------- Row 10 ----------
$ws.Range("A10").Value = ""
$ws.Range("C10").Value = 96
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "16.0"
$e10 = @'
Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .
'@
$ws.Range("E10").Value = $e10
$ws.Range("F10").Value = 0
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.00"

# ---------- Row 11 ----------
$ws.Range("C11").Value = 37
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32"
$ws.Range("E11").Value = " 50/63 A rating"
$ws.Range("F11").Value = 900
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "33300.00"

# ---------- Row 12 ----------
$ws.Range("A12").Value = ""
$ws.Range("C12").Value = 35
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "38"
$ws.Range("E12").Value = "Grand Total"
$ws.Range("F12").Value = 0
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.00"

# ---------- Row 14 ----------
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "309996.00"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "309996.00"

# ---------- Row 16 ----------
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "309996.00"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "309996.00"
